$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$ts = "2025-11-30 02:01:47"

# New listing to append (scraper re-run picked up one additional job ahead of
# the previous lowest-priority entry).
$newTitle    = "空き室情報を拾ってくスクリプト作成"
$newCategory = "システム開発"
$newPrice    = "5,000 円 ~ 10,000 円 / 固定"
$newDeadline = "期限情報なし"
$newUrl      = "https://www.lancers.jp/work/detail/5444064"
$newScore    = 10

# Previously-last row's hyperlink target, needed again once that row shifts
# down to row 13.
$shiftedUrl = "https://www.lancers.jp/work/detail/5443921"

# Insert a fresh row above the current last data row (row 12), pushing the
# existing row 12 down to row 13.
$ws.Rows.Item(12).Insert()

# Refresh the "fetched at" timestamp on every pre-existing data row (2-11),
# and on the row that just shifted down to 13.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $ts
}
$ws.Cells.Item(13, 1).Value = $ts

# Populate the newly inserted row 12 with the new listing.
$ws.Cells.Item(12, 1).Value = $ts
$ws.Cells.Item(12, 2).Value = $newTitle
$ws.Cells.Item(12, 3).Value = $newCategory
$ws.Cells.Item(12, 4).Value = $newPrice
$ws.Cells.Item(12, 5).Value = $newDeadline
$ws.Cells.Item(12, 6).Value = $newUrl
$ws.Cells.Item(12, 7).Value = $newScore

# Row 13 (the listing that shifted down) needs its hyperlink re-registered
# against its new location. Re-apply the Hyperlink cell style afterwards,
# since Hyperlinks.Add() can otherwise leave the cell on a freshly minted
# (but equivalent) style index instead of the shared one.
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), $shiftedUrl) | Out-Null
$ws.Cells.Item(13, 6).Style = "Hyperlink"

Write-Output "done"
